$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

# Insert two new rows above the old row 17 ("NumOfRetries" ...), shifting
# everything below down by two (old row 17 becomes row 19, old row 18
# becomes row 20, etc.) - this makes room for the new
# "PopUpAdsBlockerDelay" setting plus one blank spacer row, matching the
# author's intent of handling the pop-up ads blocker.
$null = $ws.Rows("17:18").Insert()

# The newly inserted rows lose their explicit 14.25pt row height; restore
# it to match the surrounding rows.
$ws.Rows("16:17").RowHeight = 14.25

# Populate the new setting row.
$ws.Range("A17").Value = "PopUpAdsBlockerDelay"
$ws.Range("B17").Value = 1000

# Update the view: select A17 and activate the sheet (drops the old
# scroll position / selection on A21).
$null = $ws.Activate()
$null = $ws.Range("A17").Select()
